$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 712, shifting existing rows 712:756 down to 716:760.
$ws.Rows("712:715").Insert()

# Helper to populate a row with the record fields (columns A..R)
function Set-Row {
    param($r, $D, $H, $I, $J, $K, $L, $M, $N, $O, $P)
    $ws.Cells.Item($r, 1).Value = 6
    $ws.Cells.Item($r, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"
    $ws.Cells.Item($r, 4).Value = $D
    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = 100112031
    $ws.Cells.Item($r, 7).Value = "Poroto verde"
    $ws.Cells.Item($r, 8).Value = $H
    $ws.Cells.Item($r, 9).Value = $I
    $ws.Cells.Item($r, 10).Value = $J
    $ws.Cells.Item($r, 11).Value = $K
    $ws.Cells.Item($r, 12).Value = $L
    $ws.Cells.Item($r, 13).Value = $M
    $ws.Cells.Item($r, 14).Value = $N
    $ws.Cells.Item($r, 15).Value = $O
    $ws.Cells.Item($r, 16).Value = $P
    $ws.Cells.Item($r, 17).Value = 25
    $ws.Cells.Item($r, 18).Value = "Hortaliza"
}

Set-Row 712 44610 "Magnum" "Primera" 3000 350 400 380 "`$/caja 25 kilos" "Perú" 15
Set-Row 713 44610 "Magnum" "Primera" 1500 15000 17000 16067 "`$/saco 25 kilos" "Región Metropolitana" 643
Set-Row 714 44610 "Magnum" "Segunda" 600 300 300 300 "`$/caja 25 kilos" "Perú" 12
Set-Row 715 44610 "Sin especificar" "Primera" 350 25000 30000 27857 "`$/malla 25 kilos" "Provincia del Elquí" 1114

# Ensure the date cells keep the date number format consistent with column D.
$ws.Range("D712:D715").NumberFormat = $ws.Range("D716").NumberFormat
